# Apply motilal_portfolio_change_engine update:
# Insert a new "Industry" column at column C (between "Stock Name" and
# "Mutual Fund"), shifting the existing Mutual Fund/Status/Jan_2026/
# Dec_2025/Oct_2025/MoM/QoQ columns one position to the right, and
# populate the new Industry column with each holding's industry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts C:I -> D:J and
# copies formatting (e.g. the header style) from the inserted-before column.
$ws.Range("C1").EntireColumn.Insert()

# Header
$ws.Range("C1").Value = "Industry"

$industries = @{
    2  = "Banks"
    3  = "Finance"
    4  = "Construction"
    5  = "Metals & Minerals Trading"
    6  = "Banks"
    7  = "Insurance"
    8  = "Power"
    9  = "Pharmaceuticals & Biotechnology"
    10 = "Automobiles"
    11 = "Banks"
    12 = "Power"
    13 = "Insurance"
    14 = "Beverages"
    15 = "Personal Products"
    16 = "Power"
    17 = "Entertainment"
    18 = "Banks"
    19 = "Ferrous Metals"
    20 = "Finance"
    21 = "Auto Components"
    22 = "Pharmaceuticals & Biotechnology"
    23 = "Finance"
    24 = "IT - Software"
    25 = "Realty"
    26 = "Healthcare Services"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}

Write-Output "Done"
